{"js": "// Fill in the delivery date (\"Data da Entrega\") blank with 09 / 06 /2020.\nconst body = context.document.body;\n\nconst target = \"Data da Entrega:         /      /2020\";\nconst filled = \"Data da Entrega:     09  / 06  /2020\";\n\nconst results = body.search(target, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(filled, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fill in the delivery date (\"Data da Entrega\") blank with 09 / 06 /2020.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Data da Entrega:         /      /2020\"\n$find.Replacement.Text = \"Data da Entrega:     09  / 06  /2020\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
